# FA80_TestData_ManuallyManageTransfers_21C.xlsx - "Add files via upload"
#
# The refreshed test-data file no longer ships real/example credentials on
# the Input_Value sheet: the URL / UserName / Password cells (S2:U2) are
# wiped out (and, since those three shared strings become unreferenced
# anywhere in the workbook, they naturally drop out of sharedStrings.xml
# too). The last thing the author did before saving was leave the cursor
# selecting that same S2:U2 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Clear the stored URL / UserName / Password test values.
$ws.Range("S2:U2").ClearContents()

# Leave the selection on the (now empty) credential cells, matching the
# selection that was active when the workbook was last saved.
$ws.Range("S2:U2").Select()
